$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7253948450088501
$ws.Range("B1").Value = 1.537975430488586
$ws.Range("C1").Value = 4.353761672973633
$ws.Range("D1").Value = 2.50617241859436
$ws.Range("E1").Value = 0.9278985261917114
